$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.193.33"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.826.69"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'236.29"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "'0.6131"
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.07107"
$ws.Range("E8").Value = "  -4.39%  "
$ws.Range("D9").Value = "'0.2815"
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("D10").Value = "'23.60"
$ws.Range("E10").Value = "  -5.33%  "
$ws.Range("D11").Value = "'0.07666"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").Value = "1.823.25"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").Value = "'4.819"
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("D14").Value = "'0.00001009"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "'0.6331"
$ws.Range("E15").Value = "  -6.23%  "
$ws.Range("D16").Value = "2.064.76"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("E17").Value = "  -3.07%  "
$ws.Range("D18").Value = "'5.872"
$ws.Range("E18").Value = "  -5.72%  "
$ws.Range("D19").Value = "29.167.09"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "'227.78"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("E21").Value = "  -4.11%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'7.012"
$ws.Range("E23").Value = "  -4.61%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'155.12"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").Value = "'0.1316"
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("D27").Value = "'8.043"
$ws.Range("E27").Value = "  -5.21%  "
$ws.Range("E28").Value = "  -4.44%  "
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").Value = "'0.06358"
$ws.Range("E30").Value = "  -11.04%  "
$ws.Range("D31").Value = "'1.452"
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").Value = "'3.823"
$ws.Range("E32").Value = "  -5.62%  "
$ws.Range("D33").Value = "'3.804"
$ws.Range("E33").Value = "  -5.62%  "
$ws.Range("D34").Value = "'1.130"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").Value = "'1.747"
$ws.Range("E35").Value = "  -4.28%  "
$ws.Range("D36").Value = "'0.6503"
$ws.Range("E36").Value = "  -6.52%  "
$ws.Range("D37").Value = "'2.546"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("D39").Value = "1.218.09"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("D40").Value = "'6.585"
$ws.Range("E40").Value = "  -3.13%  "
$ws.Range("E41").Value = "  -5.23%  "
$ws.Range("D42").Value = "'0.9190"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("D43").Value = "'0.9996"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "'101.39"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("D45").Value = "1.974.29"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").Value = "'63.05"
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").Value = "'1.627"
$ws.Range("E48").Value = "  -4.81%  "
$ws.Range("D49").Value = "'8.592"
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").Value = "'0.05520"
$ws.Range("E51").Value = "  -2.63%  "
